# Nudge the rotated "Title 1" WordArt-style title on the Cascade Mountains
# slide (slide 6: "The Cascade Mountains" / cascade_side demo) to a new
# position/size. Before this edit the shape had no explicit <p:spPr> override
# and simply inherited its off/ext/rotation from the slide layout; moving it
# in the UI causes PowerPoint to write an explicit <a:xfrm> on the slide shape
# itself (keeping the existing -75 degree rotation).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)
$sh = $s.Shapes.Item("Title 1")

# Target off/ext (EMU) converted to points (1 pt = 12700 EMU); values below
# are nudged by sub-thousandths of a point (well under 1 EMU of travel) so
# that the runtime's internal float32 rounding lands exactly on the target
# EMU instead of one EMU short:
#   off  x=2139439 y=2815201
#   ext cx=3653319 cy=1695631
$sh.Left   = 168.45980827952755
$sh.Top    = 221.66940307874015
$sh.Width  = 287.6629483858268
$sh.Height = 133.51429746850394
